$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = '8 Sep 2021'
$ws.Range("B42").Value = '0,6'
$ws.Range("C42").Value = '-7.5,6'
$ws.Range("D42").Value = '-17.5,6'
$ws.Range("E42").Value = '-40,4,8,7'
$ws.Range("F42").Value = '-27.5,2,7,5,6,4'
$ws.Range("G42").Value = '-10,3,7,5,6'
$ws.Range("H42").Value = '-37.5,4,9,7'
$ws.Range("I42").Value = '-20,3,7,5,6'
$ws.Range("J42").Value = '-27.5,3,7,6,6'
$ws.Range("A43").Value = '15 Sep 2021'
$ws.Range("B43").Value = '0,6'
$ws.Range("C43").Value = '-15,6'
$ws.Range("D43").Value = '-25,6'
$ws.Range("E43").Value = '-45,6'
$ws.Range("F43").Value = '-32.5,3,7,5,4'
$ws.Range("G43").Value = '-17.5,5,7'
$ws.Range("H43").Value = '-42.5,4,9,6'
$ws.Range("I43").Value = '-25,4,7,6'
$ws.Range("J43").Value = '-32.5,4,6,6'
$ws.Range("A44").Value = '22 Sep 2021'
$ws.Range("B44").Value = '0,6'
$ws.Range("C44").Value = '-15,6'
$ws.Range("D44").Value = '-25,6'
$ws.Range("E44").Value = '-45,6'
$ws.Range("F44").Value = '-32.5,6'
$ws.Range("G44").Value = '-17.5,6'
$ws.Range("H44").Value = '-42.5,3,5,6,5'
$ws.Range("I44").Value = '-25,5,7'
$ws.Range("J44").Value = '-32.5,5,8'
$ws.Range("A45").Value = '29 Sep 2021'
$ws.Range("B45").Value = '0,6'
$ws.Range("C45").Value = '-15,6'
$ws.Range("D45").Value = '-25,6'
$ws.Range("E45").Value = '-45,6'
$ws.Range("F45").Value = '-32.5,3,7,6,5'
$ws.Range("G45").Value = '-17.5,3,9,6,5'
$ws.Range("H45").Value = '-42.5,4,6,4'
$ws.Range("I45").Value = '-25,6'
$ws.Range("J45").Value = '-32.5,4,8,6'
$ws.Range("A46").Value = '6 Oct 2021'
$ws.Range("B46").Value = '0,6'
$ws.Range("C46").Value = '-15,6'
$ws.Range("D46").Value = '-25,6'
$ws.Range("E46").Value = '-45,6'
$ws.Range("F46").Value = '-32.5,3,5,4,4'
$ws.Range("G46").Value = '-17.5,3,8,5,4'
$ws.Range("H46").Value = '-42.5,2,8,6,5,5'
$ws.Range("I46").Value = '-25,5,5'
$ws.Range("J46").Value = '-32.5,4,6,5'
$ws.Range("A47").Value = '13 Oct 2021'
$ws.Range("B47").Value = '0,6'
$ws.Range("C47").Value = '-20,6'
$ws.Range("D47").Value = '-30,6'
$ws.Range("E47").Value = '-50,6'
$ws.Range("F47").Value = '-37.5,5,6'
$ws.Range("G47").Value = '-22.5,5,6'
$ws.Range("H47").Value = '-47.5,5,9'
$ws.Range("I47").Value = '-30,5,7'
$ws.Range("J47").Value = '-30,4,6,6'
$ws.Range("A48").Value = '20 Oct 2021'
$ws.Range("B48").Value = '0,6'
$ws.Range("C48").Value = '-20,6'
$ws.Range("D48").Value = '-27.5,6'
$ws.Range("E48").Value = '-50,6'
$ws.Range("F48").Value = '-37.5,5,5'
$ws.Range("G48").Value = '-22.5,6'
$ws.Range("H48").Value = '-47.5,6'
$ws.Range("I48").Value = '-30,5,8'
$ws.Range("J48").Value = '-30,3,7,8,6'
$ws.Range("A49").Value = '3 Nov 2021'
$ws.Range("B49").Value = '0,6'
$ws.Range("C49").Value = '-17.5,5,6'
$ws.Range("D49").Value = '-25,6'
$ws.Range("E49").Value = '-47.5,6'
$ws.Range("F49").Value = '-37.5,4,4,6'
$ws.Range("G49").Value = '-22.5,6'
$ws.Range("H49").Value = '-45,4,7,4'
$ws.Range("I49").Value = '-27.5,4,7,6'
$ws.Range("J49").Value = '-30,3,5,5,5'
$ws.Range("A50").Value = '10 Nov 2021'
$ws.Range("B50").Value = '0,6'
$ws.Range("C50").Value = '-17.5,6'
$ws.Range("D50").Value = '-25,6'
$ws.Range("E50").Value = '-47.5,6'
$ws.Range("F50").Value = '-37.5,5,9'
$ws.Range("G50").Value = '-22.5,6'
$ws.Range("H50").Value = '-45,3,9,5,5'
$ws.Range("I50").Value = '-27.5,4,8,6'
$ws.Range("J50").Value = '-30,3,9,6,5'
$ws.Range("A51").Value = '17 Nov 2021'
$ws.Range("B51").Value = '0,6'
$ws.Range("C51").Value = '-17.5,6'
$ws.Range("D51").Value = '-25,6'
$ws.Range("E51").Value = '-47.5,5,6'
$ws.Range("F51").Value = '-37.5,5,6'
$ws.Range("G51").Value = '-20,3,9,7,5'
$ws.Range("H51").Value = '-45,2,8,7,4,5'
$ws.Range("I51").Value = '-27.5,4,5,5'
$ws.Range("J51").Value = '-30,4,5,4'
$ws.Range("A52").Value = '24 Nov 2021'
$ws.Range("B52").Value = '0,6'
$ws.Range("C52").Value = '-15,5,5'
$ws.Range("D52").Value = '-22.5,6'
$ws.Range("E52").Value = '-47.5,4,9,5'
$ws.Range("F52").Value = '-37.5,5,6'
$ws.Range("G52").Value = '-20,5,9'
$ws.Range("H52").Value = '-45,3,9,5,4'
$ws.Range("I52").Value = '-27.5,4,7,6'
$ws.Range("J52").Value = '-30,4,7,6'
$ws.Range("A53").Value = '1 Dec 2021'
$ws.Range("B53").Value = '0,6'
$ws.Range("C53").Value = '-15,6'
$ws.Range("D53").Value = '-22.5,6'
$ws.Range("E53").Value = '-47.5,5,8'
$ws.Range("F53").Value = '-37.5,3,8,6,6'
$ws.Range("G53").Value = '-20,6'
$ws.Range("H53").Value = '-45,4,6,6'
$ws.Range("I53").Value = '-27.5,4,8,6'
$ws.Range("J53").Value = '-30,3,9,4,6'
$ws.Range("A54").Value = '8 Dec 2021'
$ws.Range("B54").Value = '0,6'
$ws.Range("C54").Value = '-12.5,6'
$ws.Range("D54").Value = '-20,6'
$ws.Range("E54").Value = '-47.5,6'
$ws.Range("F54").Value = '-37.5,3,9,6,5'
$ws.Range("G54").Value = '-17.5,4,7,8'
$ws.Range("H54").Value = '-45,5,5'
$ws.Range("I54").Value = '-27.5,4,7,6'
$ws.Range("J54").Value = '-30,3,7,5,7'
$ws.Range("A55").Value = '15 Dec 2021'
$ws.Range("B55").Value = '0,6'
$ws.Range("C55").Value = '-17.5,6'
$ws.Range("D55").Value = '-20,6'
$ws.Range("E55").Value = '-45,5,5'
$ws.Range("F55").Value = '-37.5,5,9'
$ws.Range("G55").Value = '-17.5,5,8'
$ws.Range("H55").Value = '-42.5,4,9,5'
$ws.Range("I55").Value = '-27.5,4,6,6'
$ws.Range("J55").Value = '-30,4,6,8'

$ws.Range("H49").Select() | Out-Null
